$d = $word.ActiveDocument
$vt = [char]11

function Replace-Text($findText, $replaceText, $label) {
    $ok = $d.Content.Find.Execute(
        $findText, $true, $false, $false, $false, $false,
        $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Warning "Replacement failed: $label"
    }
}

# 1) Update the date on the first page.
Replace-Text "Date: 2020-06-11" "Date: 2020-06-16" "date"

# 2) Remove the "other attached packages:" / "[1] dplyr_1.0.0" lines from the
#    R session info block (the blank-line run that already preceded them is
#    kept, so only the two text lines plus one of the two trailing line
#    breaks need to go).
$removeOther = "other attached packages:" + $vt + "[1] dplyr_1.0.0" + $vt + $vt
Replace-Text $removeOther "" "other attached packages block"

# 3) Rewrite the "loaded via a namespace (and not attached):" package table.
#    The first six rows are re-sorted/re-grouped text substitutions (same
#    run/style layout), while the last two rows of the old listing are
#    dropped entirely (including their leading line break).

Replace-Text `
    " [1] rstudioapi_0.11      knitr_1.28           magrittr_1.5         tidyselect_1.1.0    " `
    " [1] crayon_1.3.4         digest_0.6.25        lifecycle_0.2.0      magrittr_1.5        " `
    "namespace row 1"

Replace-Text `
    " [5] R6_2.4.1             rlang_0.4.6          stringr_1.4.0        tools_3.5.0         " `
    " [5] evaluate_0.14        pillar_1.4.4         rlang_0.4.6          stringi_1.4.6       " `
    "namespace row 2"

Replace-Text `
    " [9] xfun_0.14            tinytex_0.23         htmltools_0.4.0.9003 ellipsis_0.3.1      " `
    " [9] rstudioapi_0.11      vctrs_0.3.1          ellipsis_0.3.1       rmarkdown_2.2       " `
    "namespace row 3"

Replace-Text `
    "[13] yaml_2.2.1           digest_0.6.25        tibble_3.0.1         lifecycle_0.2.0     " `
    "[13] tools_3.5.0          stringr_1.4.0        tinytex_0.23         xfun_0.14           " `
    "namespace row 4"

Replace-Text `
    "[17] crayon_1.3.4         purrr_0.3.4          base64enc_0.1-3      vctrs_0.3.1         " `
    "[17] yaml_2.2.1           rsconnect_0.8.16     compiler_3.5.0       pkgconfig_2.0.3     " `
    "namespace row 5"

Replace-Text `
    "[21] rsconnect_0.8.16     glue_1.4.1           evaluate_0.14        rmarkdown_2.2       " `
    "[21] base64enc_0.1-3      htmltools_0.4.0.9003 knitr_1.28           tibble_3.0.1        " `
    "namespace row 6"

$removeRow7 = $vt + "[25] stringi_1.4.6        compiler_3.5.0       pillar_1.4.4         generics_0.0.2      "
Replace-Text $removeRow7 "" "namespace row 7 (removed)"

$removeRow8 = $vt + "[29] pkgconfig_2.0.3     "
Replace-Text $removeRow8 "" "namespace row 8 (removed)"
